# Update price list values in "Hoja1" (ARANDELA CHAPISTA y comunes)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date cell (A1): updated list date
$ws.Range("A1").Value = 45436

# Updated prices in column D
$ws.Range("D33").Value = 5634.136
$ws.Range("D34").Value = 4233.155
$ws.Range("D35").Value = 3652.985
$ws.Range("D36").Value = 3280.528
$ws.Range("D37").Value = 3280.528
$ws.Range("D38").Value = 2766.95
$ws.Range("D39").Value = 2766.95
$ws.Range("D40").Value = 2766.95
$ws.Range("D41").Value = 2766.95
$ws.Range("D42").Value = 2766.95
$ws.Range("D43").Value = 2766.95
$ws.Range("D44").Value = 2766.95
$ws.Range("D45").Value = 3302.011
$ws.Range("D46").Value = 3302.011
$ws.Range("D47").Value = 3302.011

$ws.Range("D51").Value = 4813.344
$ws.Range("D52").Value = 4376.426
$ws.Range("D53").Value = 3652.985
$ws.Range("D54").Value = 3652.985
